# Insert a new weekly price record as row 20 ("Madrigal" / "Primera", 2021-10-08,
# Provincia del Elqui) into the Alcachofa price sheet. This shifts all existing
# records previously in rows 20-60 down to rows 21-61.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 20, pushing rows 20:60 down to 21:61.
$ws.Rows("20:20").Insert()

# Populate the newly inserted row 20 with the new record's data.
$ws.Range("A20").Value = 5
$ws.Range("B20").Value = "Macroferia Regional de Talca"
$ws.Range("C20").Value = "Maule"
$ws.Range("D20").Value = 44477
$ws.Range("E20").Value = 7
$ws.Range("F20").Value = 100112013
$ws.Range("G20").Value = "Alcachofa"
$ws.Range("H20").Value = "Madrigal"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 300
$ws.Range("K20").Value = 10000
$ws.Range("L20").Value = 10000
$ws.Range("M20").Value = 10000
$ws.Range("N20").Value = "`$/caja 40 unidades"
$ws.Range("O20").Value = "Provincia del Elquí"
$ws.Range("P20").Value = 250
$ws.Range("Q20").Value = 40
$ws.Range("R20").Value = "Hortaliza"
